# Edit kbase.xlsx per commit: "Add files via upload"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (EM30001 / 3.0.3 firmware issue): update body text, remove phoneticPr (handled by engine), adjust row height
$ws.Range("D8").Value = "iCUP 시스템을 3.0.3으로 업데이트한 경우 일부 고객은 아래와 같은 문제를 경험할 수 있습니다.`n<br>`n<br>1. 후방 카메라가 나오지 않음`n<br>2. 예기치 않게 중앙화면이 재시동되거나 꺼짐`n<br>3. 차량의 시동을 끄기 전 에어컨을 껐으나, 시동 후 자동으로 켜짐 등`n<br>`n<br>이 문제를 해결하기 위해 6월 셋째 주 개선된 펌웨어가 제공되며, 서비스센터를 통해 우선 업데이트 제공 예정임을 안내합니다.`n<br>OTA(Over-The-Air) 업데이트의 경우 7월 첫째 주 배포 예정임을 안내 합니다.`n<br>`n<br>이 문제를 임시로 해결하기 위해 중앙 화면의 홈 버튼을 20초 이상 길게 눌러 Volvo 로고가 나오면 떼도록 안내합니다. (강제 리셋)"
$ws.Rows.Item(8).RowHeight = 174.75

# --- Row 26 (Privacy Mode notice): prepend model applicability header, adjust row height
$ws.Range("D26").Value = "<br>`n<br><strong><h1>해당하는 모델:</strong> C40, 2022-2024년식/EC40, 2025년식 및 이후 연식/EX40, 2025년 모델 이후/S60, S60 Recharge, 2023년식 및 이후 연식/S90, S90 Recharge, 2022년식 및 이후 연식`n<br>V60, V60 Recharge, V60 Cross Country, 2023년식 및 이후 연식/V90, V90 Recharge, V90 Cross Country, 2022년식 및 이후 연식/XC40Recharge Pure Electric, 2021-2024년식`n<br>XC40 (Mild Hybrid), 2023년식 및 이후 연식/XC60, XC60 Recharge, 2022년식 및 이후 연식/XC90, XC90 Recharge, 2023년식 및 이후 연식</h1>`n<br>`n<br>고객이 Volvo Cars앱에 로그인 후 차량과 연결했으나, Privacy Mode에 대해 문의하는 경우 아래 단계를 수행합니다.`n<br>`n<br><img src=`"https://www.volvocars.com/images/support/imgcaf2b6bd5242fff1c0a801520a655252_1_--_--_VOICEpnghigh.png`" width=`"335`" height=`"725`">`n<br>1. 차량의 중앙 화면 우측 하단 설정을 누릅니다.`n<br>2. 개인정보 보호로 이동합니다.`n<br>3. Volvo 개인 정보 설정으로 이동합니다.`n<br>4. Volvo Cars앱을 허용(스위치를 오른쪽으로 이동) 합니다.`n<br><em><strong>스위치가 이동하지 않는 문제나 회색으로 보이는 경우 고객이 차량에 Volvo ID로 로그인을 했는지 확인합니다.</em></strong>"
$ws.Rows.Item(26).RowHeight = 261.75

# --- Row 33 (Sensus Navigation 2016 update guide): add images, close header tags, add new "current map version" section
$ws.Range("D33").Value = "Sensus Navigation(2016)이 구동 중인 경우 업데이트 방법에 대해 고객이 보고할 수 있습니다.`n<br>`n<br>차량이 네트워크가 동작하는 경우, 무선(OTA)로 업데이트할 수 있습니다.`n<br><img src=`"https://www.volvocars.com/images/support/imgd92f5d3ee39dababc0a8015252766f2b_1_--_--_VOICEpnghigh.png?imwidth=1600`">`n<br><h3>인터넷을 통해 업데이트 방법</h3>`n<br>1. 소프트웨어 업데이트는 운전 중이거나 차량의 시동이 켜져 있는 경우 차량에 직접 다운로드됩니다.`n<br>2. 업데이트가 다운로드되면 업데이트를 설치할 준비가 되었다는 알림을 받게 됩니다. '지금 설치'를 클릭하여 설치를 시작합니다.`n<br>즉시 설치하거나 원하는 시간을 선택하여 설치를 시작할 수 있습니다. 설치를 완료하는 데 약 90분이 소요됩니다.`n<br>3.소프트웨가 최신 상태입니다. 메시지가 중앙 화면에 나타나면 소프트웨어 업데이트가 완료된 것입니다.`n<br><h4>차량에서 직접 다운로드하기</h4>`n<br>1. <img src=`"https://www.volvocars.com/images/support/imgc0eb4a9734129c4bc0a801522aaf58ce_2_--_--_VOICEpnghigh.png`" width=`"64`" height=`"64`">`n<br>앱 보기에서 다운로드 센터 버튼을 누릅니다.`n<br>2. 지도 버튼을 누릅니다.`n<br>3.설치 버튼을 누른 후 확인을 선택합니다.`n<br>4. 선택한 지도 업데이트의 설치가 시작됩니다.`n<br>`n<br><h4>진행 중인 지도 다운로드가 취소될 경우`n<br>지도 다운로드 중에 차량 시동을 끄면, 차량을 다시 시동하여 인터넷에 다시 연결할 때 지도 다운로드가 다시 시작됩니다.`n<br><h3>지도를 USB로 설치하는 방법</h3>`n<br>새 지도는 인터넷에 연결된 컴퓨터에서 USB 메모리로 다운로드한 후 USB 메모리에서 차량의 내비게이션 시스템으로 전송할 수 있습니다.`n<br>`n<br><img src=`"https://www.volvocars.com/images/support/imgac2b29e3d3ad1130c0a801525f9e3341_1_--_--_VOICEpnghigh.png`">`n<br>USB 메모리`n<br>USB 메모리로 업데이트를 관리해야 하는 경우에 다음 요구 사항이 적용됩니다.`n<br>USB 표준:`t최소 2.0`n<br>파일 시스템: FAT32, exFAT 또는 NFTS`n<br>용량: 최대 128GB`n<br>`n<br>지도 파일은 아래 링크에서 다운로드 받을 수 있습니다.`n<br><img src=`"https://az685612.vo.msecnd.net/swdlimages/SPA_KR_v01_644x386.png`">`n<br><a href=`"https://www.volvocars.com/kr/support/downloads/maps/spa/daehanmingug`" target=`"_blank`">Sensus Navigation 2016 지도 다운로드</a>`n<br>`n<br><h3>현재 지도 버전 확인</h3>`n<br>1. <img src=`"https://www.volvocars.com/images/support/imgc0eb4a9734129c4bc0a801522aaf58ce_2_--_--_VOICEpnghigh.png`" width=`"64`" height=`"64`">`n<br>앱 보기에서 다운로드 센터 버튼을 누릅니다.`n<br>2. 지도 버튼을 누릅니다.`n<br>> 대한민국을 누릅니다.`n<br>3.<img=`"https://www.volvocars.com/images/support/imgf13a34931b25bc6ac0a8015256d334bc_1_--_--_VOICEpnghigh.png`" width=`"64`" height=`"64`">`n<br>화살표를 터치하여 선택한 지역을 확대합니다.`n<br>> 상세 지도 정보가 표시됩니다.`n<br>4. 지도 이미지 아래에서 지도 정보를 확인합니다."
$ws.Range("E33").Value = "All"
$ws.Rows.Item(33).RowHeight = 409.6

# --- New row 34: Bluetooth media sound issue
$ws.Range("A33:E33").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Range("A34").Value = "Bluetooth를 통해 미디어 사운드가 재생되지 않는 경우"
$ws.Range("B34").Value2 = 45454
$ws.Range("C34").Value = "HT203027"
$ws.Range("D34").Value = "<br>지원 모델: iCup 지원 차량`n<br>`n<br>고객은 차량에 블루투스로 연결하면 음악이나 기타 미디어 사운드가 재생되지 않는다고 보고할 수 있습니다.`n<br>홈 보기 또는 앱 보기<img src=`"https://www.volvocars.com/images/support/img7a8807df5cacedebc0a801525a46ef6f_2_--_--_VOICEpnghigh.png`" width=`"64`" height=`"64`">에서 블루투스 미디어 플레이어 앱<img src=`"https://www.volvocars.com/images/support/img963609385e11fb97c0a801525c7b6fab_3_--_--_VOICEpnghigh.png`" width=`"64`" height=`"64`">을 시작합니다.`n<br> 앱 사용 시 홈 보기에서 단축 명령을 사용하여 앱을 제어할 수도 있습니다.`n<br>`n<br><em>참고`n<br>Bluetooth를 통해 전화에서 미디어를 스트리밍하려면 먼저 블루투스 미디어 플레이어를 시작해야 합니다.</em>`n<br>`n<br>음성 컨트롤을 사용하여 미디어 플레이어를 제어할 수도 있습니다."
$ws.Range("E34").Value = "All"
$ws.Rows.Item(34).RowHeight = 209.25

# --- Update selection / active cell to match final view state
$ws.Range("B8").Select()
